$d = $word.ActiveDocument

# --- Target: first paragraph of the document (the "**ID__...__ID**" placeholder line) ---
$p1 = $d.Paragraphs(1)

# 1) Remove the trailing run that contains only a single space character
#    (it sits right before the paragraph mark at the end of paragraph 1).
$pEnd = $p1.Range.End
$trailingSpace = $d.Range($pEnd - 2, $pEnd - 1)
if ($trailingSpace.Text -eq " ") {
    $trailingSpace.Delete()
}

# 2) Update the placeholder id text in the (now single) remaining run.
$d.Content.Find.Execute("**ID__AFFARS_5317_topic_4__ID**", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "**ID__AFFARS_5317_106__ID**", 2)

# 3) Adjust paragraph indentation: 120 twips -> 225 twips (1 pt = 20 twips).
$p1 = $d.Paragraphs(1)
$p1.Format.LeftIndent = 225 / 20

# 4) Add paragraph border spacing (w:space="5") on all four sides, matching
#    the border-spacing-only <w:pBdr> already used elsewhere in the document.
$borders = $p1.Format.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
